# Auto-generated Excel COM-interop script
# Applies numeric cell updates across multiple worksheets per the commit diff.
$wb = $excel.ActiveWorkbook

# Sheet ALC, row 17 (Leve Item ID 38956)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2186.25
$ws.Range("I17").Value = 1499.3334
$ws.Range("J17").Value = 2598.4
$ws.Range("K17").Value = 4498.0002
$ws.Range("L17").Value = 7795.200000000001
$ws.Range("M17").Value = -4330.0002
$ws.Range("N17").Value = -8131.200000000001

# Sheet ALC, row 28 (Leve Item ID 27772)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 648.05884
$ws.Range("I28").Value = 594.8125
$ws.Range("K28").Value = 594.8125
$ws.Range("M28").Value = -109.8125

# Sheet ALC, row 31 (Leve Item ID 4576)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 1852.2858
$ws.Range("I31").Value = 1852.2858
$ws.Range("K31").Value = 5556.857400000001
$ws.Range("M31").Value = -5326.857400000001

# Sheet ALC, row 99 (Leve Item ID 19883)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 429.75
$ws.Range("I99").Value = 429.75
$ws.Range("K99").Value = 1289.25
$ws.Range("M99").Value = 208.75

# Sheet ALC, row 121 (Leve Item ID 39731)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 298
$ws.Range("I121").Value = 298
$ws.Range("K121").Value = 894
$ws.Range("M121").Value = 853

# Sheet ALC, row 127 (Leve Item ID 36114)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 2800.348
$ws.Range("I127").Value = 2700.3635
$ws.Range("J127").Value = 5000
$ws.Range("K127").Value = 8101.0905
$ws.Range("L127").Value = 15000
$ws.Range("M127").Value = -3141.0905
$ws.Range("N127").Value = -24920

# Sheet ALC, row 131 (Leve Item ID 36108)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 13982.917
$ws.Range("J131").Value = 76500
$ws.Range("L131").Value = 229500
$ws.Range("N131").Value = -239580

# Sheet ALC, row 132 (Leve Item ID 44049)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 10177.679
$ws.Range("I132").Value = 10349.038
$ws.Range("K132").Value = 31047.114
$ws.Range("M132").Value = -28517.114

# Sheet ALC, row 137 (Leve Item ID 44013)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 23414.143
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()

# Sheet ALC, row 138 (Leve Item ID 44169)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1085.7916
$ws.Range("I138").Value = 1085.7916
$ws.Range("K138").Value = 3257.3748
$ws.Range("M138").Value = 1882.6252

# Sheet ALC, row 141 (Leve Item ID 44161)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 5616.4
$ws.Range("I141").Value = 5725.4116
$ws.Range("K141").Value = 17176.2348
$ws.Range("M141").Value = -11996.2348

# Sheet ARM, row 4 (Leve Item ID 5071)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 171
$ws.Range("I4").Value = 171.8
$ws.Range("K4").Value = 171.8
$ws.Range("M4").Value = -55.80000000000001

# Sheet ARM, row 32 (Leve Item ID 44147)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 292799.03
$ws.Range("I32").Value = 309822.88
$ws.Range("K32").Value = 309822.88
$ws.Range("M32").Value = -309535.88

# Sheet ARM, row 45 (Leve Item ID 27714)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3324.6191
$ws.Range("I45").Value = 3070.125
$ws.Range("K45").Value = 3070.125
$ws.Range("M45").Value = -2693.125

# Sheet ARM, row 61 (Leve Item ID 43999)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3354.7805
$ws.Range("I61").Value = 3098.389
$ws.Range("K61").Value = 3098.389
$ws.Range("M61").Value = -2886.389

# Sheet ARM, row 132 (Leve Item ID 43997)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 596878.4399999999
$ws.Range("I132").Value = 611329.1
$ws.Range("K132").Value = 1833987.3
$ws.Range("M132").Value = -1831457.3

# Sheet ARM, row 136 (Leve Item ID 43999)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3354.7805
$ws.Range("I136").Value = 3098.389
$ws.Range("K136").Value = 9295.167000000001
$ws.Range("M136").Value = -6745.167000000001

# Sheet BSM, row 94 (Leve Item ID 19939)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2968
$ws.Range("J94").Value = 734.25
$ws.Range("L94").Value = 734.25
$ws.Range("N94").Value = -1636.25

# Sheet BSM, row 105 (Leve Item ID 19947)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2507.425
$ws.Range("I105").Value = 1748.1111
$ws.Range("K105").Value = 1748.1111
$ws.Range("M105").Value = -1.111100000000079

# Sheet BSM, row 107 (Leve Item ID 27706)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1483.1428
$ws.Range("I107").Value = 1465
$ws.Range("K107").Value = 1465
$ws.Range("M107").Value = 455

# Sheet CRP, row 7 (Leve Item ID 5361)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1915.2142
$ws.Range("I7").Value = 46.166668
$ws.Range("K7").Value = 46.166668
$ws.Range("M7").Value = 66.833332

# Sheet CRP, row 16 (Leve Item ID 27691)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 202599.6
$ws.Range("I16").Value = 2999.3333
$ws.Range("K16").Value = 2999.3333
$ws.Range("M16").Value = -2712.3333

# Sheet CRP, row 22 (Leve Item ID 5367)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1414.4762
$ws.Range("J22").Value = 2132.6667
$ws.Range("L22").Value = 2132.6667
$ws.Range("N22").Value = -2832.6667

# Sheet CRP, row 31 (Leve Item ID 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2133.8928
$ws.Range("I31").Value = 2161.074
$ws.Range("K31").Value = 2161.074
$ws.Range("M31").Value = -1866.074

# Sheet CRP, row 34 (Leve Item ID 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2133.8928
$ws.Range("I34").Value = 2161.074
$ws.Range("K34").Value = 2161.074
$ws.Range("M34").Value = -1959.074

# Sheet CRP, row 69 (Leve Item ID 11911)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H69").Value = 13500
$ws.Range("I69").Value = 13500
$ws.Range("K69").Value = 13500
$ws.Range("M69").Value = -12751

# Sheet CRP, row 72 (Leve Item ID 11911)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H72").Value = 13500
$ws.Range("I72").Value = 13500
$ws.Range("K72").Value = 40500
$ws.Range("M72").Value = -36756

# Sheet CRP, row 93 (Leve Item ID 19516)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H93").Value = 9935.666999999999
$ws.Range("I93").Value = 9935.666999999999
$ws.Range("K93").Value = 9935.666999999999
$ws.Range("M93").Value = -8063.666999999999

# Sheet CRP, row 105 (Leve Item ID 19928)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 23517.4
$ws.Range("I105").Value = 33863.332
$ws.Range("J105").Value = 7998.5
$ws.Range("K105").Value = 33863.332
$ws.Range("L105").Value = 7998.5
$ws.Range("M105").Value = -32116.332
$ws.Range("N105").Value = -11492.5

# Sheet CRP, row 113 (Leve Item ID 27691)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 202599.6
$ws.Range("I113").Value = 2999.3333
$ws.Range("K113").Value = 2999.3333
$ws.Range("M113").Value = -829.3332999999998

# Sheet CRP, row 114 (Leve Item ID 27112)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H114").Value = 79999.5
$ws.Range("J114").Value = 79999.5
$ws.Range("L114").Value = 79999.5
$ws.Range("N114").Value = -88677.5

# Sheet CRP, row 122 (Leve Item ID 36196)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 8504.76
$ws.Range("I122").Value = 1524
$ws.Range("K122").Value = 4572
$ws.Range("M122").Value = -2122

# Sheet CUL, row 61 (Leve Item ID 4727)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H61").Value = 35.333332
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 35.333332
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 105.999996
$ws.Range("N61").Value = -535.999996
$ws.Range("M61").ClearContents()

# Sheet CUL, row 121 (Leve Item ID 27878)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 4298.2144
$ws.Range("I121").Value = 566
$ws.Range("K121").Value = 1698
$ws.Range("M121").Value = -388

# Sheet GSM, row 25 (Leve Item ID 2687)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H25").Value = 2336
$ws.Range("I25").Value = 1004
$ws.Range("K25").Value = 1004
$ws.Range("M25").Value = -475

# Sheet GSM, row 97 (Leve Item ID 19940)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1796
$ws.Range("I97").Value = 1852.2858
$ws.Range("J97").Value = 1533.3334
$ws.Range("K97").Value = 1852.2858
$ws.Range("L97").Value = 1533.3334
$ws.Range("M97").Value = -1356.2858
$ws.Range("N97").Value = -2525.3334

# Sheet GSM, row 113 (Leve Item ID 27710)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2224.25
$ws.Range("I113").Value = 2224.25
$ws.Range("K113").Value = 2224.25
$ws.Range("M113").Value = -54.25

# Sheet GSM, row 122 (Leve Item ID 36182)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 53580.95
$ws.Range("I122").Value = 80278.53999999999
$ws.Range("J122").Value = 3999.7144
$ws.Range("K122").Value = 240835.62
$ws.Range("L122").Value = 11999.1432
$ws.Range("M122").Value = -238385.62
$ws.Range("N122").Value = -16899.1432

# Sheet GSM, row 132 (Leve Item ID 44008)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 15479.8125
$ws.Range("I132").Value = 17130.535
$ws.Range("J132").Value = 3924.75
$ws.Range("K132").Value = 51391.605
$ws.Range("L132").Value = 11774.25
$ws.Range("M132").Value = -48861.605
$ws.Range("N132").Value = -16834.25

# Sheet LTW, row 20 (Leve Item ID 4308)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 40
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

# Sheet LTW, row 40 (Leve Item ID 36248)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3182.1428
$ws.Range("I40").Value = 2506
$ws.Range("J40").Value = 4872.5
$ws.Range("K40").Value = 2506
$ws.Range("L40").Value = 4872.5
$ws.Range("M40").Value = -2370
$ws.Range("N40").Value = -5144.5

# Sheet LTW, row 46 (Leve Item ID 5282)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4538.75
$ws.Range("I46").Value = 1097.5
$ws.Range("J46").Value = 4921.1113
$ws.Range("K46").Value = 1097.5
$ws.Range("L46").Value = 4921.1113
$ws.Range("M46").Value = -909.5
$ws.Range("N46").Value = -5297.1113

# Sheet LTW, row 132 (Leve Item ID 44058)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2780251.8
$ws.Range("I132").Value = 6251590
$ws.Range("K132").Value = 18754770
$ws.Range("M132").Value = -18752240

# Sheet WVR, row 107 (Leve Item ID 27746)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1677.1724
$ws.Range("I107").Value = 1012.9583
$ws.Range("K107").Value = 3038.8749
$ws.Range("M107").Value = -1118.8749

# Sheet WVR, row 123 (Leve Item ID 34127)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 59999
$ws.Range("J123").Value = 59999
$ws.Range("L123").Value = 59999
$ws.Range("N123").Value = -69799

# Sheet WVR, row 132 (Leve Item ID 44029)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2666
$ws.Range("I132").Value = 2666
$ws.Range("K132").Value = 7998
$ws.Range("M132").Value = -5468
